$d = $word.ActiveDocument

function Replace-WithBreak($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# Paragraph 1: Bibliografia / "Programa" reading list - split run-on sentences with manual line breaks
Replace-WithBreak "1999.GIESECKE" "1999.^lGIESECKE"
Replace-WithBreak "2002.RIBEIRO" "2002.^lRIBEIRO"
Replace-WithBreak "2013.SILVA" "2013.^lSILVA"
Replace-WithBreak "2013.CRUZ" "2013.^lCRUZ"
Replace-WithBreak "2010.LIMA" "2010.^lLIMA"
Replace-WithBreak "2015.LEAKE" "2015.^lLEAKE"
Replace-WithBreak "2010.FISCHER" "2010.^lFISCHER"
Replace-WithBreak "2011.PROVENZA" "2011.^lPROVENZA"
Replace-WithBreak "1991.PROVENZA" "1991.^lPROVENZA"

# Paragraph 2: "Critério" answer text
Replace-WithBreak "independência.Para a formação" "independência.^lPara a formação"

# Paragraph 3: "Norma de recuperação" answer text
Replace-WithBreak "deficiências.- Definição" "deficiências.^l- Definição"
Replace-WithBreak "importantes.- Pesquisa" "importantes.^l- Pesquisa"
Replace-WithBreak "usuários.- Avaliação" "usuários.^l- Avaliação"
Replace-WithBreak "recebido.- Refinamento" "recebido.^l- Refinamento"
Replace-WithBreak "usuários.- Implementação" "usuários.^l- Implementação"
